$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these numeric-looking
# strings as text (matches the source inlineStr cells), same as
# typing  '314.90  into a cell in the Excel UI.

$ws.Range("D2").Value = "'28.378.83"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "'1.822.67"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'314.90"
$ws.Range("E5").Value = "  -1.20%  "
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Value = "'0.5137"
$ws.Range("E7").Value = "  -3.67%  "
$ws.Range("D8").Value = "'0.3923"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("D9").Value = "'0.07651"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.107"
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'41.58"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "'20.97"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").Value = "'1.002"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D15").Value = "'7.514"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "'1.824.03"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "'93.31"
$ws.Range("E17").Value = "  +3.64%  "
$ws.Range("D18").Value = "'0.00001095"
$ws.Range("E18").Value = "  +1.87%  "
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("D20").Value = "'17.68"
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").Value = "'6.139"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "'28.387.71"
$ws.Range("E23").Value = "  -0.64%  "
$ws.Range("D24").Value = "'11.20"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("D25").Value = "'2.256"
$ws.Range("E25").Value = "  +7.35%  "
$ws.Range("D26").Value = "'20.79"
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("D27").Value = "'156.13"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "'2.031.22"
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "'2.386"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").Value = "'124.19"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").Value = "'1.109"
$ws.Range("E31").Value = "  -2.02%  "
$ws.Range("D32").Value = "'0.1087"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "'5.650"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").Value = "'3.657"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'0.07080"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "'0.2201"
$ws.Range("E36").Value = "  -3.06%  "
$ws.Range("D37").Value = "'0.02319"
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.160"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'8.782"
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("D40").Value = "'0.6239"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("D41").Value = "'11.21"
$ws.Range("E41").Value = "  -1.81%  "
$ws.Range("D42").Value = "'1.171"
$ws.Range("E42").Value = "  -2.59%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").Value = "'13.30"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "'0.5874"
$ws.Range("E46").Value = "  -0.08%  "
$ws.Range("D47").Value = "'3.706"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").Value = "'124.83"
$ws.Range("E48").Value = "  -1.11%  "
$ws.Range("D49").Value = "'1.980"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "'1.198"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'0.06898"
$ws.Range("E51").Value = "  -0.49%  "
